# Update countries & provincias Spain
# - Re-sorted a few tied rows (same "Casos totales") whose underlying data
#   changed slightly, causing their relative order (and therefore which
#   country name sits on which row) to swap.
# - Refreshed a handful of numeric statistics for Chequia, Georgia and
#   Guinea Ecuatorial.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------
# 1) Swap the two rows for "Belice" / "Santa Lucia" (rows 188-189)
# ---------------------------------------------------------------------
$ws.Cells.Item(188, 1).Value = "Santa Lucia"
$ws.Cells.Item(188, 2).Value = 18
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 15
$ws.Cells.Item(188, 5).Value = 3
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

$ws.Cells.Item(189, 1).Value = "Belice"
$ws.Cells.Item(189, 2).Value = 18
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 13
$ws.Cells.Item(189, 5).Value = 3
$ws.Cells.Item(189, 6).Value = 1
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 2

# ---------------------------------------------------------------------
# 2) Swap the two rows for "San Vicente y las Granadinas" / "Namibia"
#    (rows 194-195). Underlying numbers are identical for these two,
#    only the country names trade places.
# ---------------------------------------------------------------------
$ws.Cells.Item(194, 1).Value = "Namibia"
$ws.Cells.Item(194, 2).Value = 16
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 8
$ws.Cells.Item(194, 5).Value = 8
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

$ws.Cells.Item(195, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(195, 2).Value = 16
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 8
$ws.Cells.Item(195, 5).Value = 8
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

# ---------------------------------------------------------------------
# 3) Swap the two rows for "San Cristobal y Nieves" / "Burundi"
#    (rows 198-199)
# ---------------------------------------------------------------------
$ws.Cells.Item(198, 1).Value = "Burundi"
$ws.Cells.Item(198, 2).Value = 15
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 7
$ws.Cells.Item(198, 5).Value = 7
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 1

$ws.Cells.Item(199, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(199, 2).Value = 15
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 8
$ws.Cells.Item(199, 5).Value = 7
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0

# ---------------------------------------------------------------------
# 4) Refresh plain numeric values (no reordering involved)
# ---------------------------------------------------------------------
# Chequia - row 47
$ws.Cells.Item(47, 4).Value = 3592
$ws.Cells.Item(47, 5).Value = 3940
$ws.Cells.Item(47, 6).Value = 58
$ws.Cells.Item(47, 7).Value = 1
$ws.Cells.Item(47, 8).Value = 249

# Georgia - row 111
$ws.Cells.Item(111, 2).Value = 593
$ws.Cells.Item(111, 3).Value = 4
$ws.Cells.Item(111, 5).Value = 361

# Guinea Ecuatorial - row 130
$ws.Cells.Item(130, 4).Value = 13
$ws.Cells.Item(130, 5).Value = 301
